$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The edit reshuffles/renumbers the existing Q&A rows (one row of old content
# was dropped), so every destination cell's new text already exists verbatim
# somewhere else in the original sheet. We reproduce that reshuffle with
# Copy / PasteSpecial(xlPasteValues) so Excel never has to "guess" the type of
# a freshly-typed literal (this matters for the "01/01/2020" text, which Excel
# would otherwise auto-convert to a real date if typed in directly).
# The copy order below is topologically sorted so a cell is always read before
# it gets overwritten by a later step.
# ---------------------------------------------------------------------------
$ws.Range("B26").Copy()
$ws.Range("B24").PasteSpecial(-4163)
$ws.Range("C26").Copy()
$ws.Range("C24").PasteSpecial(-4163)
$ws.Range("B14").Copy()
$ws.Range("B18").PasteSpecial(-4163)
$ws.Range("A16").Copy()
$ws.Range("A14").PasteSpecial(-4163)
$ws.Range("B13").Copy()
$ws.Range("B10").PasteSpecial(-4163)
$ws.Range("C14").Copy()
$ws.Range("C18").PasteSpecial(-4163)
$ws.Range("C25").Copy()
$ws.Range("C23").PasteSpecial(-4163)
$ws.Range("B25").Copy()
$ws.Range("B23").PasteSpecial(-4163)
$ws.Range("C13").Copy()
$ws.Range("C10").PasteSpecial(-4163)
$ws.Range("B20").Copy()
$ws.Range("B19").PasteSpecial(-4163)
$ws.Range("C13").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$ws.Range("C20").Copy()
$ws.Range("C19").PasteSpecial(-4163)
$ws.Range("A15").Copy()
$ws.Range("A13").PasteSpecial(-4163)
$ws.Range("B13").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("A18").Copy()
$ws.Range("A16").PasteSpecial(-4163)
$ws.Range("B21").Copy()
$ws.Range("B20").PasteSpecial(-4163)
$ws.Range("C8").Copy()
$ws.Range("C13").PasteSpecial(-4163)
$ws.Range("C21").Copy()
$ws.Range("C20").PasteSpecial(-4163)
$ws.Range("A17").Copy()
$ws.Range("A15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("B13").PasteSpecial(-4163)
$ws.Range("A20").Copy()
$ws.Range("A18").PasteSpecial(-4163)
$ws.Range("B22").Copy()
$ws.Range("B21").PasteSpecial(-4163)
$ws.Range("C22").Copy()
$ws.Range("C21").PasteSpecial(-4163)
$ws.Range("A19").Copy()
$ws.Range("A17").PasteSpecial(-4163)
$ws.Range("A22").Copy()
$ws.Range("A20").PasteSpecial(-4163)
$ws.Range("A21").Copy()
$ws.Range("A19").PasteSpecial(-4163)
$ws.Range("A24").Copy()
$ws.Range("A22").PasteSpecial(-4163)
$ws.Range("A23").Copy()
$ws.Range("A21").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# Cells that held content before but must end up empty in the new layout
$ws.Range("B14").ClearContents()
$ws.Range("C14").ClearContents()
$ws.Range("B17").ClearContents()
$ws.Range("C17").ClearContents()
$ws.Range("B22").ClearContents()
$ws.Range("C22").ClearContents()
$ws.Range("A23").ClearContents()
$ws.Range("A24").ClearContents()

# Drop the two now-unused trailing rows (sheet shrinks from 26 to 24 rows)
$ws.Rows.Item(25).Delete()
$ws.Rows.Item(25).Delete()

# Row heights for the reshuffled rows (10 keeps its height; 13-24 all change)
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(14).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(16).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(19).RowHeight = 60
$ws.Rows.Item(20).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).RowHeight = 30
$ws.Rows.Item(24).RowHeight = 30
